# Auto-generated Excel COM-interop script applying the commit diff
# to Sheets/Alpha_Profits.xlsx (workbook sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 70.22221999999999
$ws.Range("I5").Value = 75.25
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 75.25
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = 39.75
$ws.Range("N5").Value = -260
$ws.Range("H106").Value = 9329.833000000001
$ws.Range("J106").Value = 9998
$ws.Range("L106").Value = 9998
$ws.Range("N106").Value = -11260
$ws.Range("H112").Value = 3357.1428
$ws.Range("J112").Value = 3420
$ws.Range("L112").Value = 10260
$ws.Range("N112").Value = -12476
$ws.Range("H125").Value = 7638.25
$ws.Range("J125").Value = 6933.3335
$ws.Range("L125").Value = 62400.0015
$ws.Range("N125").Value = -67320.0015
$ws.Range("H137").Value = 90002180
$ws.Range("I137").Value = 150001730
$ws.Range("J137").Value = 2850
$ws.Range("K137").Value = 450005190
$ws.Range("L137").Value = 8550
$ws.Range("M137").Value = -450002640
$ws.Range("N137").Value = -13650
$ws.Range("H138").Value = 2842.6333
$ws.Range("J138").Value = 2929.7222
$ws.Range("L138").Value = 8789.1666
$ws.Range("N138").Value = -19069.1666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 3005
$ws.Range("I25").Value = 3005
$ws.Range("K25").Value = 3005
$ws.Range("M25").Value = -2603
$ws.Range("H32").Value = 3128.5806
$ws.Range("I32").Value = 3207.8276
$ws.Range("J32").Value = 1979.5
$ws.Range("K32").Value = 3207.8276
$ws.Range("L32").Value = 1979.5
$ws.Range("M32").Value = -2920.8276
$ws.Range("N32").Value = -2553.5
$ws.Range("H95").Value = 23803.2
$ws.Range("J95").Value = 23803.2
$ws.Range("L95").Value = 23803.2
$ws.Range("N95").Value = -29295.2
$ws.Range("H101").Value = 15301
$ws.Range("J101").Value = 15301
$ws.Range("L101").Value = 15301
$ws.Range("N101").Value = -21791
$ws.Range("H122").Value = 1615.3334
$ws.Range("I122").Value = 1019.25
$ws.Range("K122").Value = 3057.75
$ws.Range("M122").Value = -607.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 5879.5
$ws.Range("I24").Value = 1166.6666
$ws.Range("K24").Value = 1166.6666
$ws.Range("M24").Value = -931.6666
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H34").Value = 15000
$ws.Range("I34").Value = 15000
$ws.Range("K34").Value = 15000
$ws.Range("M34").Value = -14886
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H81").Value = 7499
$ws.Range("J81").Value = 7499
$ws.Range("L81").Value = 7499
$ws.Range("N81").Value = -9621
$ws.Range("H84").Value = 7499
$ws.Range("J84").Value = 7499
$ws.Range("L84").Value = 22497
$ws.Range("N84").Value = -33105
$ws.Range("H86").Value = 2279.5417
$ws.Range("I86").Value = 2122.25
$ws.Range("J86").Value = 2594.125
$ws.Range("K86").Value = 2122.25
$ws.Range("L86").Value = 2594.125
$ws.Range("M86").Value = -999.25
$ws.Range("N86").Value = -4840.125
$ws.Range("H89").Value = 2279.5417
$ws.Range("I89").Value = 2122.25
$ws.Range("J89").Value = 2594.125
$ws.Range("K89").Value = 10611.25
$ws.Range("L89").Value = 12970.625
$ws.Range("M89").Value = -4995.25
$ws.Range("N89").Value = -24202.625
$ws.Range("H110").Value = 39955.332
$ws.Range("J110").Value = 39955.332
$ws.Range("L110").Value = 39955.332
$ws.Range("N110").Value = -48135.332
$ws.Range("H134").Value = 14708351
$ws.Range("I134").Value = 20835596
$ws.Range("J134").Value = 2962.4
$ws.Range("K134").Value = 62506788
$ws.Range("L134").Value = 8887.200000000001
$ws.Range("M134").Value = -62504253
$ws.Range("N134").Value = -13957.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2416.6
$ws.Range("J31").Value = 3216.6667
$ws.Range("L31").Value = 3216.6667
$ws.Range("N31").Value = -3806.6667
$ws.Range("H34").Value = 2416.6
$ws.Range("J34").Value = 3216.6667
$ws.Range("L34").Value = 3216.6667
$ws.Range("N34").Value = -3620.6667
$ws.Range("H54").Value = 7561.364
$ws.Range("I54").Value = 7308.3
$ws.Range("J54").Value = 10092
$ws.Range("K54").Value = 7308.3
$ws.Range("L54").Value = 10092
$ws.Range("M54").Value = -6650.3
$ws.Range("N54").Value = -11408
$ws.Range("H109").Value = 39998.5
$ws.Range("J109").Value = 39998.5
$ws.Range("L109").Value = 39998.5
$ws.Range("N109").Value = -42078.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 145.26414
$ws.Range("I38").Value = 26
$ws.Range("J38").Value = 160.48936
$ws.Range("K38").Value = 78
$ws.Range("L38").Value = 481.46808
$ws.Range("M38").Value = 269
$ws.Range("N38").Value = -1175.46808
$ws.Range("H88").Value = 3099
$ws.Range("H91").Value = 3099
$ws.Range("H113").Value = 576.8
$ws.Range("I113").Value = 235
$ws.Range("J113").Value = 804.6667
$ws.Range("K113").Value = 705
$ws.Range("L113").Value = 2414.0001
$ws.Range("M113").Value = 1465
$ws.Range("N113").Value = -6754.0001
$ws.Range("H132").Value = 1698.25
$ws.Range("I132").Value = 1698
$ws.Range("J132").Value = 1699
$ws.Range("K132").Value = 15282
$ws.Range("L132").Value = 15291
$ws.Range("M132").Value = -12752
$ws.Range("N132").Value = -20351
$ws.Range("H141").Value = 3800
$ws.Range("I141").Value = 4666.6665
$ws.Range("J141").Value = 1200
$ws.Range("K141").Value = 13999.9995
$ws.Range("L141").Value = 3600
$ws.Range("M141").Value = -8819.999500000002
$ws.Range("N141").Value = -13960

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3529.8235
$ws.Range("I122").Value = 3660.9333
$ws.Range("K122").Value = 10982.7999
$ws.Range("M122").Value = -8532.7999
$ws.Range("H132").Value = 3121.75
$ws.Range("I132").Value = 2749.5
$ws.Range("J132").Value = 3494
$ws.Range("K132").Value = 8248.5
$ws.Range("L132").Value = 10482
$ws.Range("M132").Value = -5718.5
$ws.Range("N132").Value = -15542

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4579.4
$ws.Range("I7").Value = 5999
$ws.Range("K7").Value = 5999
$ws.Range("M7").Value = -5887
$ws.Range("H16").Value = 796.3333
$ws.Range("J16").Value = 990
$ws.Range("L16").Value = 990
$ws.Range("N16").Value = -1330
$ws.Range("H22").Value = 4808.6
$ws.Range("I22").Value = 3010.625
$ws.Range("J22").Value = 6007.25
$ws.Range("K22").Value = 3010.625
$ws.Range("L22").Value = 6007.25
$ws.Range("M22").Value = -2715.625
$ws.Range("N22").Value = -6597.25
$ws.Range("H27").Value = 4808.6
$ws.Range("I27").Value = 3010.625
$ws.Range("J27").Value = 6007.25
$ws.Range("K27").Value = 3010.625
$ws.Range("L27").Value = 6007.25
$ws.Range("M27").Value = -2903.625
$ws.Range("N27").Value = -6221.25
$ws.Range("H40").Value = 2780.889
$ws.Range("I40").Value = 2487.1428
$ws.Range("K40").Value = 2487.1428
$ws.Range("M40").Value = -2351.1428
$ws.Range("H68").Value = 4694.0586
$ws.Range("I68").Value = 2587.7144
$ws.Range("K68").Value = 2587.7144
$ws.Range("M68").Value = -1838.7144
$ws.Range("H71").Value = 4694.0586
$ws.Range("I71").Value = 2587.7144
$ws.Range("K71").Value = 12938.572
$ws.Range("M71").Value = -9194.572
$ws.Range("H82").Value = 4200.25
$ws.Range("I82").Value = 2999.8
$ws.Range("J82").Value = 6201
$ws.Range("K82").Value = 2999.8
$ws.Range("L82").Value = 6201
$ws.Range("M82").Value = -2638.8
$ws.Range("N82").Value = -6923
$ws.Range("H85").Value = 4200.25
$ws.Range("I85").Value = 2999.8
$ws.Range("J85").Value = 6201
$ws.Range("K85").Value = 2999.8
$ws.Range("L85").Value = 6201
$ws.Range("M85").Value = -1751.8
$ws.Range("N85").Value = -8697
$ws.Range("H122").Value = 3622.125
$ws.Range("I122").Value = 3496.8572
$ws.Range("K122").Value = 10490.5716
$ws.Range("M122").Value = -8040.571599999999
$ws.Range("H126").Value = 4579.4
$ws.Range("I126").Value = 5999
$ws.Range("K126").Value = 17997
$ws.Range("M126").Value = -15527

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H122").Value = 1938.0769
$ws.Range("I122").Value = 1874.5834
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 5623.7502
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -3173.7502
$ws.Range("N122").Value = -13000
$ws.Range("H132").Value = 2285.8572
$ws.Range("I132").Value = 2460.4
$ws.Range("J132").Value = 1849.5
$ws.Range("K132").Value = 7381.200000000001
$ws.Range("L132").Value = 5548.5
$ws.Range("M132").Value = -4851.200000000001
$ws.Range("N132").Value = -10608.5
